$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text columns (Coin name, Link) - safe to assign directly
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"

# Numeric-looking text columns (Price, Volume) - force text format so Excel
# does not coerce values like "1.00" or "28.20" into numbers, losing the
# literal text representation from the source diff. Reset the number format
# back to the default afterwards so no stray style survives on the cell.
$deCells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($cellRef in $deCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "68.392.70"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "2.644.17"
$ws.Range("E3").Value = "  +0.93%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "598.61"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("D6").Value = "154.52"
$ws.Range("E6").Value = "  +1.25%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.547"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "2.642.92"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").Value = "0.143"
$ws.Range("E10").Value = "  +7.47%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "5.27"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("D13").Value = "0.355"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").Value = "28.20"
$ws.Range("E14").Value = "  +2.40%  "
$ws.Range("D15").Value = "0.0000191"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("D16").Value = "3.124.39"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").Value = "68.328.64"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").Value = "2.649.63"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").Value = "11.46"
$ws.Range("E19").Value = "  +2.33%  "
$ws.Range("D20").Value = "365.62"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").Value = "7.49"
$ws.Range("E21").Value = "  +13.11%  "
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  +3.66%  "
$ws.Range("D23").Value = "4.90"
$ws.Range("E23").Value = "  +2.15%  "
$ws.Range("D24").Value = "2.08"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "73.88"
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").Value = "9.90"
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "0.0000107"
$ws.Range("E28").Value = "  +2.77%  "
$ws.Range("D29").Value = "2.779.64"
$ws.Range("E29").Value = "  +1.17%  "
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "577.13"
$ws.Range("E31").Value = "  -2.99%  "
$ws.Range("D32").Value = "8.21"
$ws.Range("E32").Value = "  +5.29%  "
$ws.Range("D33").Value = "1.43"
$ws.Range("E33").Value = "  +3.49%  "
$ws.Range("D34").Value = "1.87"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("E35").Value = "  +4.76%  "
$ws.Range("D36").Value = "1.60"
$ws.Range("E36").Value = "  +5.83%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("D38").Value = "158.91"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "19.44"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").Value = "1.90"
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("D41").Value = "0.374"
$ws.Range("E41").Value = "  +1.86%  "
$ws.Range("D42").Value = "5.44"
$ws.Range("E42").Value = "  +3.19%  "
$ws.Range("D43").Value = "0.0₆0344"
$ws.Range("E43").Value = "  +15.44%  "
$ws.Range("D44").Value = "2.68"
$ws.Range("E44").Value = "  -1.32%  "
$ws.Range("E45").Value = "  +3.52%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "40.55"
$ws.Range("E47").Value = "  +0.41%  "
$ws.Range("D48").Value = "157.51"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "3.78"
$ws.Range("E49").Value = "  +2.79%  "
$ws.Range("D50").Value = "1.72"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "22.00"
$ws.Range("E51").Value = "  +3.44%  "

foreach ($cellRef in $deCells) {
    $ws.Range($cellRef).Style = "Normal"
}
